$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 1
$ws.Range("A3").Value = 4
$ws.Range("A4").Value = 3
$ws.Range("A5").Value = 2
